$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (correlated variable names)
$ws.Range("B1").Value = "IT.NET.USER.P2:BMU"
$ws.Range("C1").Value = "SP.URB.TOTL:BMU"

# Row labels
$ws.Range("A2").Value = "SP.POP.TOTL:BMU:cor-value"
$ws.Range("A3").Value = "SP.POP.TOTL:BMU:p-value"
$ws.Range("A4").Value = "SP.URB.TOTL:BMU:cor-value"
$ws.Range("A5").Value = "SP.URB.TOTL:BMU:p-value"

# Values
$ws.Range("B2").Value = 0.8066759342260912
$ws.Range("C2").Value = 1
$ws.Range("B3").Value = 0.00048995054476078
$ws.Range("C3").Value = 0
$ws.Range("B4").Value = 0.8066759342260912
$ws.Range("B5").Value = 0.00048995054476078

# Build the bold / centered / boxed format once on a scratch cell, then
# stamp it onto every label cell via copy/paste-special so the style
# table only ever gains the one extra cell format (matches the target
# workbook's single additional cellXfs entry instead of one per
# incremental property mutation).
$scratch = $ws.Range("Z100")
$scratch.Value = "x"
$scratch.Font.Bold = $true
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160
$scratch.Borders.LineStyle = 1

$scratch.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$ws.Range("A2:A5").PasteSpecial(-4122)
$scratch.Clear()
